# Linked front-end and back-end
# Adds a new student (Jared Rathbun) to the "students" sheet and his
# per-course rows to the "class_data" sheet, then updates the active
# sheet/selection to reflect the end state captured by the diff.

$wb = $excel.ActiveWorkbook

$classData = $wb.Worksheets.Item("class_data")
$students  = $wb.Worksheets.Item("students")

# ---------------------------------------------------------------------
# class_data: append rows 9-15 for student 789456 (Jared Rathbun)
# ---------------------------------------------------------------------
$newClassRows = @(
    @(789456, "UNDG", "Day - Science & Eng ", "Problem Solving w/ Java",   "CSC1610", "B", "FA", 2019),
    @(789456, "UNDG", "Day - Science & Eng ", "Data Structures",           "CSC2820", "A", "SP", 2020),
    @(789456, "UNDG", "Day - Science & Eng ", "Object Oriented Design",    "CSC2620", "C", "FA", 2020),
    @(789456, "UNDG", "Day - Science & Eng ", "Analysis of Algorithms",    "CSC2710", "D", "SP", 2021),
    @(789456, "UNDG", "Day - Science & Eng ", "Network Security",          "CSC5055", "B", "FA", 2022),
    @(789456, "UNDG", "Day - Science & Eng ", "Web Development",           "CSC3222", "A", "SP", 2022),
    @(789456, "UNDG", "Day - Science & Eng ", "Computer Science Capstone", "CSC3333", "F", "FA", 2022)
)

$row = 9
foreach ($rec in $newClassRows) {
    $classData.Cells.Item($row, 1).Value = $rec[0]
    $classData.Cells.Item($row, 2).Value = $rec[1]
    $classData.Cells.Item($row, 3).Value = $rec[2]
    $classData.Cells.Item($row, 4).Value = $rec[3]
    $classData.Cells.Item($row, 5).Value = $rec[4]
    $classData.Cells.Item($row, 6).Value = $rec[5]
    $classData.Cells.Item($row, 7).Value = $rec[6]
    $classData.Cells.Item($row, 8).Value = $rec[7]
    $row++
}

# ---------------------------------------------------------------------
# students: append row 4 for Jared Rathbun
# ---------------------------------------------------------------------
$students.Cells.Item(4, 1).Value  = 789456
$students.Cells.Item(4, 2).Value  = "Rathbun"
$students.Cells.Item(4, 3).Value  = "Jared"
$students.Cells.Item(4, 4).Value  = "Computer Science"
$students.Cells.Item(4, 7).Value  = "Software Engineering"
$students.Cells.Item(4, 10).Value = "Cyber Security"
$students.Cells.Item(4, 13).Value = "MA"
$students.Cells.Item(4, 14).Value = "USA"
$students.Cells.Item(4, 17).Value = 7
$students.Cells.Item(4, 18).Value = 1.2
$students.Cells.Item(4, 19).Value = 2.1
$students.Cells.Item(4, 20).Value = 3.99
$students.Cells.Item(4, 21).Value = 400
$students.Cells.Item(4, 22).Value = 2
$students.Cells.Item(4, 23).Value = "N"
$students.Cells.Item(4, 24).Value = "White"

# ---------------------------------------------------------------------
# Update selections / active sheet to match the saved view state
# ---------------------------------------------------------------------
$classData.Range("D17").Select()
$students.Range("Q19").Select()
$students.Activate()

# Best-effort: mirror the saved window geometry (no-op if unsupported
# by this headless host).
try {
    $win = $excel.ActiveWindow
    $win.Left   = -108
    $win.Top    = -108
    $win.Width  = 30936
    $win.Height = 16776
} catch {
}
